# Apply the "merge data columns from multiple data frames" edit:
# - Append 11 new rows (24-34) of Date/Link data to Sheet1
# - Remove the stray empty E4 cell style marker (handled naturally by not writing it)
# - Re-style column E (and header E3) to use the "left aligned" style (style index 3)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New link rows: column D = date (Excel serial number), column E = link text
$newRows = @(
    @{ Row = 24; Date = 44881; Link = "https://zoom.us/rec/play/9rqDVSztxqxx7-GqyKcIREg9doHHfcbPMxhmBDXxY4nMt5OQcfjDJnOJXl1vyycvXAmuWPTpYCMumy3H.S780cc3gEc40E39v?continueMode=true&_x_zm_rtaid=wwSi-ScIRa2lw1QLNM3sOw.1669538040676.1ddd9f75194d27c04e9ba5973bde60bd&_x_zm_rhtaid=55" },
    @{ Row = 25; Date = 44883; Link = "https://zoom.us/rec/play/N6NVnKqOwe9i8DoJ4rm-yA5aylNS4qMFdWws485gQZwyVTG-KpSJ6AriGcMJ1tQCNDETmXHyHZi_zEZd.tMBaisORe59gaxdj?continueMode=true&_x_zm_rtaid=wwSi-ScIRa2lw1QLNM3sOw.1669538040676.1ddd9f75194d27c04e9ba5973bde60bd&_x_zm_rhtaid=55" },
    @{ Row = 26; Date = 44886; Link = "https://zoom.us/rec/play/3sSzhbXVUxDmrZiee-D-EjM01Wq3M0dWx59K8FrDZoqkv_zQq8YJtI0Y0biHe902suF-m4WfjEe_EyxJ.12UAZfkyH5YeQJXh?continueMode=true&_x_zm_rtaid=wwSi-ScIRa2lw1QLNM3sOw.1669538040676.1ddd9f75194d27c04e9ba5973bde60bd&_x_zm_rhtaid=55" },
    @{ Row = 27; Date = 44888; Link = "https://zoom.us/rec/share/2LjjCHY9RWe7FCyZwYw6RV33UcA2isYTxViuTZSmRSE_Xz01v6fmd9IXGTmaArFP.A-nETf4bqQJaTzY4" },
    @{ Row = 28; Date = 44890; Link = "https://zoom.us/rec/share/-EIvE6woZ6rGdxM3kp8S--yjsipOuKCXmq0l8JIhdmSLwU2JqspYCSUhIWrYWHzi.On8k-WEOislNXqse" },
    @{ Row = 29; Date = 44893; Link = "https://zoom.us/rec/play/l4g6LCEmLrR_3rqSp8zcxzhVHjas0Skxm-r8CUIIMs40HqfESw6oZqpsQkrKXXquJ19IEYvSSfgYdfCM.Lo2BMtRTYRKM_sjb?continueMode=true&_x_zm_rtaid=LS47qZHbS0Csx8jjRCN8dw.1670731521042.b3d473194791c9e526111e4bb6fde81d&_x_zm_rhtaid=184" },
    @{ Row = 30; Date = 44895; Link = "https://zoom.us/rec/play/kq4kxIHxTiLImTAXQTO2tXjvRDdFJdbyZgwIpQfKZFGaR7z-tS3kyDzX74OJwC9GHOZGBhgL3qfr5CWI.3tzz4IHhXSlRV5KR?continueMode=true&_x_zm_rtaid=LS47qZHbS0Csx8jjRCN8dw.1670731521042.b3d473194791c9e526111e4bb6fde81d&_x_zm_rhtaid=184" },
    @{ Row = 31; Date = 44897; Link = "https://zoom.us/rec/play/La9gRWETlob4q0PaROzEXHYBVFHR-PBqzwn8V53yYjN9vELZj-BRw2c6mb5EAPuoqp1MglBlwB9I3MPh.uJYcvUXWI2REjjL2?continueMode=true&_x_zm_rtaid=LS47qZHbS0Csx8jjRCN8dw.1670731521042.b3d473194791c9e526111e4bb6fde81d&_x_zm_rhtaid=184" },
    @{ Row = 32; Date = 44900; Link = "https://zoom.us/rec/play/qyDcgMrywA0CHSDPN5NMnXLSCWBuwFcTFl9H2gLdhN11rFC0tV3zzUOgJYivBSyMihR-M7J4ZRCtAGGM.dawIYKXdhkPzQ_L6?continueMode=true&_x_zm_rtaid=LS47qZHbS0Csx8jjRCN8dw.1670731521042.b3d473194791c9e526111e4bb6fde81d&_x_zm_rhtaid=184" },
    @{ Row = 33; Date = 44902; Link = "https://zoom.us/rec/play/CWxosEUV0uL2gVVPINeyzTKj-9eum8WwU0f9ZSFEiIvxNcNEDaC_OHbrVkZZrH-XTwOWOAgoWAvAGW9s.J7agLk-oR748wYhX?continueMode=true&_x_zm_rtaid=LS47qZHbS0Csx8jjRCN8dw.1670731521042.b3d473194791c9e526111e4bb6fde81d&_x_zm_rhtaid=184" },
    @{ Row = 34; Date = 44904; Link = "https://zoom.us/rec/play/1nf8LBaxQW3edx05w0EET8hdNkyhn4Y2o-ZxsQlDSNxbV1BnmXscVmY5lT_rxpTcdeuvGR3i4dAggYXV.tzWuwy7yvm3Qi-P5?continueMode=true&_x_zm_rtaid=LS47qZHbS0Csx8jjRCN8dw.1670731521042.b3d473194791c9e526111e4bb6fde81d&_x_zm_rhtaid=184" }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 4).Value = $r.Date
    $ws.Cells.Item($r.Row, 5).Value = $r.Link
}

# Remove the now-stray empty E4 style-only cell entirely
$ws.Range("E4").Clear()

# Column E (data + header) moves to the "left" aligned style used by data rows
$ws.Range("E3:E34").HorizontalAlignment = -4131  # xlLeft

# Update selection to mirror the final cursor position from the diff
$ws.Range("E34").Select()
